$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): relabel the Sum Sq / Pct Sum Sq / Mean Sq columns ---
$ws.Range("C2").Value2 = "eta_sq"
$ws.Range("D2").Value2 = "lower CI"
$ws.Range("E2").Value2 = "Upper CI"

# --- Data rows 3:10 -- replace Sum Sq / Pct Sum Sq(formula) / Mean Sq values
#     with eta_sq / lower CI / Upper CI values, each shown as a percentage ---

$etaSq   = @{3=0.0627511625; 4=0.0501918619; 5=0.0869342783; 6=0.1351966468; 7=0.2322911401; 8=0.0003583943; 9=0.0007860463; 10=0.0034645367}
$lowerCI = @{3=0.05718131; 4=0.04514611; 5=0.08055662; 6=0.1277144; 7=0.2237571; 8=0.00004431111; 9=0.0002486499; 10=0.002174011}
$upperCI = @{3=0.0685234141; 4=0.0554587643; 5=0.0934806045; 6=0.1427862175; 7=0.2686264561; 8=0.0009731077; 9=0.00162338; 10=0.005050481}

foreach ($r in 3..10) {
    $c = $ws.Cells.Item($r, 3)
    $c.Value2 = $etaSq[$r]
    $c.NumberFormat = "0.0%"

    $d = $ws.Cells.Item($r, 4)
    $d.Value2 = $lowerCI[$r]
    $d.NumberFormat = "0.0%"

    $e = $ws.Cells.Item($r, 5)
    $e.Value2 = $upperCI[$r]
    $e.NumberFormat = "0.0%"
}

# --- F4:F7 picked up a distinct (but visually identical) "no fill" style in the
#     accepted-publication pass -- reproduce that explicit formatting touch ---
foreach ($r in 4..7) {
    $f = $ws.Cells.Item($r, 6)
    $f.Interior.ColorIndex = 0
}

# --- Drop the Residuals / Total rows (11:12) -- eta_sq/CI replace the old
#     Sum-Sq-based Pct Sum Sq column so the grand-total row is no longer needed ---
$ws.Rows("11:12").Delete()

# --- Restore selection to match the trimmed table ---
$ws.Range("D8").Select()
